# Add temperature sensor(1,2)+widget FAN
# Updates the "Typography" sheet (new wildcard range for the Small typography)
# and the "Translation" sheet (new/updated translation rows for the
# power/sound temperature sensors and a new FAN widget text, plus two new
# "GB-DIRECTION"/"GB-ALIGNMENT" header columns used by the new row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Typography sheet: give the "Small" typography (row 6) the same
# "0-9" wildcard range already used by the other rows that reuse it.
# ---------------------------------------------------------------------
$typography = $wb.Worksheets.Item("Typography")
$typography.Range("H6").Value = "0-9"

# ---------------------------------------------------------------------
# Translation sheet
# ---------------------------------------------------------------------
$translation = $wb.Worksheets.Item("Translation")

# New header columns for the extra GB-direction/alignment info.
$translation.Range("G3").Value = "GB-DIRECTION"
$translation.Range("H3").Value = "GB-ALIGNMENT"

# Row 29: reword the power-sensor text (add a space after the colon).
$translation.Range("F29").Value = "Power: <> C°"

# Row 30: new "Sound" temperature-sensor row (text id reused from the old
# row 31), right aligned, with the new GB-direction/GB-alignment values.
$translation.Range("B30").Value = "SingleUseId38"
$translation.Range("C30").Value = "Small"
$translation.Range("D30").Value = "Right"
$translation.Range("E30").Value = "LTR"
$translation.Range("F30").Value = "Sound: <> C°"
$translation.Range("G30").Value = "LTR"
$translation.Range("H30").Value = "Right"

# Row 31: formerly row 32's "00" placeholder, now refined to "00.0".
$translation.Range("B31").Value = "SingleUseId39"
$translation.Range("C31").Value = "Small"
$translation.Range("D31").Value = "Left"
$translation.Range("E31").Value = "LTR"
$translation.Range("F31").Value = "00.0"

# Row 32: shifted up from the old row 33 (unchanged values).
$translation.Range("B32").Value = "SingleUseId40"
$translation.Range("C32").Value = "Default"
$translation.Range("D32").Value = "Left"
$translation.Range("E32").Value = "LTR"
$translation.Range("F32").Value = "P"

# Row 33: shifted up from the old row 34 (unchanged values).
$translation.Range("B33").Value = "SingleUseId41"
$translation.Range("C33").Value = "Default"
$translation.Range("D33").Value = "Left"
$translation.Range("E33").Value = "LTR"
$translation.Range("F33").Value = "S"

# Row 34: shifted up from the old row 35 (unchanged values).
$translation.Range("B34").Value = "SingleUseId42"
$translation.Range("C34").Value = "Typ2"
$translation.Range("D34").Value = "Left"
$translation.Range("E34").Value = "LTR"
$translation.Range("F34").Value = "Player"

# Row 35: brand new row for the FAN widget's "00.0" text.
$translation.Range("B35").Value = "SingleUseId43"
$translation.Range("C35").Value = "Small"
$translation.Range("D35").Value = "Left"
$translation.Range("E35").Value = "LTR"
$translation.Range("F35").Value = "00.0"
